$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 185
$wsExpo.Range("F5").Value = 3406
$wsExpo.Range("F6").Value = 347
$wsExpo.Range("F7").Value = 18
$wsExpo.Range("F8").Value = 424

# Sheet "全部类型" (All types) - same underlying rows, update column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 185
$wsAll.Range("F5").Value = 3406
$wsAll.Range("F6").Value = 347
$wsAll.Range("F9").Value = 18
$wsAll.Range("F10").Value = 424
